$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.788523
$ws.Range("H2").Value = 38.365569
$ws.Range("I2").Value = 0.09894253826004661
$ws.Range("J2").Value = 0.09894253826004659
$ws.Range("M2").Value = 21.85308466666666
$ws.Range("N2").Value = 65.559254
$ws.Range("O2").Value = 0.407053040353553
$ws.Range("P2").Value = 0.407053040353553
$ws.Range("Q2").Value = 279.4686758806139
$ws.Range("R2").Value = 2515.218082925526
$ws.Range("S2").Value = 0.04027486101904972
$ws.Range("T2").Value = 0.04027486101904971
$ws.Range("G3").Value = 12.788523
$ws.Range("H3").Value = 38.365569
$ws.Range("I3").Value = 0.09894253826004661
$ws.Range("J3").Value = 0.09894253826004659
$ws.Range("O3").Value = 0.1342711086924142
$ws.Range("P3").Value = 0.1342711086924142
$ws.Range("Q3").Value = 92.185944423111
$ws.Range("R3").Value = 829.6734998079991
$ws.Range("S3").Value = 0.01328512430901807
$ws.Range("T3").Value = 0.01328512430901807
$ws.Range("G4").Value = 12.788523
$ws.Range("H4").Value = 38.365569
$ws.Range("I4").Value = 0.09894253826004661
$ws.Range("J4").Value = 0.09894253826004659
$ws.Range("M4").Value = 11.375406
$ws.Range("N4").Value = 34.126218
$ws.Range("O4").Value = 0.2118874139822907
$ws.Range("P4").Value = 0.2118874139822907
$ws.Range("Q4").Value = 145.474641265338
$ws.Range("R4").Value = 1309.271771388042
$ws.Range("S4").Value = 0.02096467856476513
$ws.Range("T4").Value = 0.02096467856476513
$ws.Range("G5").Value = 12.788523
$ws.Range("H5").Value = 38.365569
$ws.Range("I5").Value = 0.09894253826004661
$ws.Range("J5").Value = 0.09894253826004659
$ws.Range("M5").Value = 3.401340666666667
$ws.Range("N5").Value = 10.204022
$ws.Range("O5").Value = 0.06335609277882483
$ws.Range("P5").Value = 0.06335609277882483
$ws.Range("Q5").Value = 43.498123346502
$ws.Range("R5").Value = 391.483110118518
$ws.Range("S5").Value = 0.006268612633775938
$ws.Range("T5").Value = 0.006268612633775937
$ws.Range("G6").Value = 12.788523
$ws.Range("H6").Value = 38.365569
$ws.Range("I6").Value = 0.09894253826004661
$ws.Range("J6").Value = 0.09894253826004659
$ws.Range("M6").Value = 9.847764666666666
$ws.Range("N6").Value = 29.543294
$ws.Range("O6").Value = 0.1834323441929172
$ws.Range("P6").Value = 0.1834323441929172
$ws.Range("Q6").Value = 125.938364938254
$ws.Range("R6").Value = 1133.445284444286
$ws.Range("S6").Value = 0.01814926173343774
$ws.Range("T6").Value = 0.01814926173343774
$ws.Range("I7").Value = 0.1332855577638913
$ws.Range("J7").Value = 0.1332855577638912
$ws.Range("M7").Value = 21.85308466666666
$ws.Range("N7").Value = 65.559254
$ws.Range("O7").Value = 0.407053040353553
$ws.Range("P7").Value = 0.407053040353553
$ws.Range("Q7").Value = 376.4724353885424
$ws.Range("R7").Value = 3388.251918496881
$ws.Range("S7").Value = 0.05425429152301105
$ws.Range("T7").Value = 0.05425429152301103
$ws.Range("I8").Value = 0.1332855577638913
$ws.Range("J8").Value = 0.1332855577638912
$ws.Range("O8").Value = 0.1342711086924142
$ws.Range("P8").Value = 0.1342711086924142
$ws.Range("S8").Value = 0.0178963996136445
$ws.Range("T8").Value = 0.01789639961364449
$ws.Range("I9").Value = 0.1332855577638913
$ws.Range("J9").Value = 0.1332855577638912
$ws.Range("M9").Value = 11.375406
$ws.Range("N9").Value = 34.126218
$ws.Range("O9").Value = 0.2118874139822907
$ws.Range("P9").Value = 0.2118874139822907
$ws.Range("Q9").Value = 195.968984043966
$ws.Range("R9").Value = 1763.720856395694
$ws.Range("S9").Value = 0.02824153215577815
$ws.Range("T9").Value = 0.02824153215577814
$ws.Range("I10").Value = 0.1332855577638913
$ws.Range("J10").Value = 0.1332855577638912
$ws.Range("M10").Value = 3.401340666666667
$ws.Range("N10").Value = 10.204022
$ws.Range("O10").Value = 0.06335609277882483
$ws.Range("P10").Value = 0.06335609277882483
$ws.Range("Q10").Value = 58.59635030469178
$ws.Range("R10").Value = 527.367152742226
$ws.Range("S10").Value = 0.008444452163766512
$ws.Range("T10").Value = 0.008444452163766509
$ws.Range("I11").Value = 0.1332855577638913
$ws.Range("J11").Value = 0.1332855577638912
$ws.Range("M11").Value = 9.847764666666666
$ws.Range("N11").Value = 29.543294
$ws.Range("O11").Value = 0.1834323441929172
$ws.Range("P11").Value = 0.1834323441929172
$ws.Range("Q11").Value = 169.6516534733558
$ws.Range("R11").Value = 1526.864881260202
$ws.Range("S11").Value = 0.02444888230769105
$ws.Range("T11").Value = 0.02444888230769104
$ws.Range("G12").Value = 53.92730466666666
$ws.Range("H12").Value = 161.781914
$ws.Range("I12").Value = 0.4172260084485798
$ws.Range("J12").Value = 0.4172260084485797
$ws.Range("M12").Value = 21.85308466666666
$ws.Range("N12").Value = 65.559254
$ws.Range("O12").Value = 0.407053040353553
$ws.Range("P12").Value = 0.407053040353553
$ws.Range("Q12").Value = 1178.477954725795
$ws.Range("R12").Value = 10606.30159253216
$ws.Range("S12").Value = 0.1698331152535716
$ws.Range("T12").Value = 0.1698331152535715
$ws.Range("G13").Value = 53.92730466666666
$ws.Range("H13").Value = 161.781914
$ws.Range("I13").Value = 0.4172260084485798
$ws.Range("J13").Value = 0.4172260084485797
$ws.Range("O13").Value = 0.1342711086924142
$ws.Range("P13").Value = 0.1342711086924142
$ws.Range("Q13").Value = 388.7344543923882
$ws.Range("R13").Value = 3498.610089531494
$ws.Range("S13").Value = 0.05602139872970138
$ws.Range("T13").Value = 0.05602139872970137
$ws.Range("G14").Value = 53.92730466666666
$ws.Range("H14").Value = 161.781914
$ws.Range("I14").Value = 0.4172260084485798
$ws.Range("J14").Value = 0.4172260084485797
$ws.Range("M14").Value = 11.375406
$ws.Range("N14").Value = 34.126218
$ws.Range("O14").Value = 0.2118874139822907
$ws.Range("P14").Value = 0.2118874139822907
$ws.Range("Q14").Value = 613.444985069028
$ws.Range("R14").Value = 5521.004865621252
$ws.Range("S14").Value = 0.08840493997632294
$ws.Range("T14").Value = 0.08840493997632291
$ws.Range("G15").Value = 53.92730466666666
$ws.Range("H15").Value = 161.781914
$ws.Range("I15").Value = 0.4172260084485798
$ws.Range("J15").Value = 0.4172260084485797
$ws.Range("M15").Value = 3.401340666666667
$ws.Range("N15").Value = 10.204022
$ws.Range("O15").Value = 0.06335609277882483
$ws.Range("P15").Value = 0.06335609277882483
$ws.Range("Q15").Value = 183.4251344064564
$ws.Range("R15").Value = 1650.826209658108
$ws.Range("S15").Value = 0.02643380970100697
$ws.Range("T15").Value = 0.02643380970100696
$ws.Range("G16").Value = 53.92730466666666
$ws.Range("H16").Value = 161.781914
$ws.Range("I16").Value = 0.4172260084485798
$ws.Range("J16").Value = 0.4172260084485797
$ws.Range("M16").Value = 9.847764666666666
$ws.Range("N16").Value = 29.543294
$ws.Range("O16").Value = 0.1834323441929172
$ws.Range("P16").Value = 0.1834323441929172
$ws.Range("Q16").Value = 531.0634054649685
$ws.Range("R16").Value = 4779.570649184716
$ws.Range("S16").Value = 0.07653274478797684
$ws.Range("T16").Value = 0.07653274478797682
$ws.Range("G17").Value = 3.523547333333334
$ws.Range("H17").Value = 10.570642
$ws.Range("I17").Value = 0.0272610618786406
$ws.Range("J17").Value = 0.0272610618786406
$ws.Range("M17").Value = 21.85308466666666
$ws.Range("N17").Value = 65.559254
$ws.Range("O17").Value = 0.407053040353553
$ws.Range("P17").Value = 0.407053040353553
$ws.Range("Q17").Value = 77.00037820234088
$ws.Range("R17").Value = 693.003403821068
$ws.Range("S17").Value = 0.011096698120967
$ws.Range("T17").Value = 0.011096698120967
$ws.Range("G18").Value = 3.523547333333334
$ws.Range("H18").Value = 10.570642
$ws.Range("I18").Value = 0.0272610618786406
$ws.Range("J18").Value = 0.0272610618786406
$ws.Range("O18").Value = 0.1342711086924142
$ws.Range("P18").Value = 0.1342711086924142
$ws.Range("Q18").Value = 25.39945689137578
$ws.Range("R18").Value = 228.595112022382
$ws.Range("S18").Value = 0.003660373002577582
$ws.Range("T18").Value = 0.003660373002577581
$ws.Range("G19").Value = 3.523547333333334
$ws.Range("H19").Value = 10.570642
$ws.Range("I19").Value = 0.0272610618786406
$ws.Range("J19").Value = 0.0272610618786406
$ws.Range("M19").Value = 11.375406
$ws.Range("N19").Value = 34.126218
$ws.Range("O19").Value = 0.2118874139822907
$ws.Range("P19").Value = 0.2118874139822907
$ws.Range("Q19").Value = 40.081781476884
$ws.Range("R19").Value = 360.7360332919561
$ws.Range("S19").Value = 0.005776275903876365
$ws.Range("T19").Value = 0.005776275903876363
$ws.Range("G20").Value = 3.523547333333334
$ws.Range("H20").Value = 10.570642
$ws.Range("I20").Value = 0.0272610618786406
$ws.Range("J20").Value = 0.0272610618786406
$ws.Range("M20").Value = 3.401340666666667
$ws.Range("N20").Value = 10.204022
$ws.Range("O20").Value = 0.06335609277882483
$ws.Range("P20").Value = 0.06335609277882483
$ws.Range("Q20").Value = 11.98478483579156
$ws.Range("R20").Value = 107.863063522124
$ws.Range("S20").Value = 0.001727154365632438
$ws.Range("T20").Value = 0.001727154365632438
$ws.Range("G21").Value = 3.523547333333334
$ws.Range("H21").Value = 10.570642
$ws.Range("I21").Value = 0.0272610618786406
$ws.Range("J21").Value = 0.0272610618786406
$ws.Range("M21").Value = 9.847764666666666
$ws.Range("N21").Value = 29.543294
$ws.Range("O21").Value = 0.1834323441929172
$ws.Range("P21").Value = 0.1834323441929172
$ws.Range("Q21").Value = 34.69906493052756
$ws.Range("R21").Value = 312.291584374748
$ws.Range("S21").Value = 0.005000560485587216
$ws.Range("T21").Value = 0.005000560485587215
$ws.Range("G22").Value = 41.78521799999999
$ws.Range("H22").Value = 125.355654
$ws.Range("I22").Value = 0.3232848336488418
$ws.Range("J22").Value = 0.3232848336488418
$ws.Range("M22").Value = 21.85308466666666
$ws.Range("N22").Value = 65.559254
$ws.Range("O22").Value = 0.407053040353553
$ws.Range("P22").Value = 0.407053040353553
$ws.Range("Q22").Value = 913.1359067691237
$ws.Range("R22").Value = 8218.223160922114
$ws.Range("S22").Value = 0.1315940744369537
$ws.Range("T22").Value = 0.1315940744369536
$ws.Range("G23").Value = 41.78521799999999
$ws.Range("H23").Value = 125.355654
$ws.Range("I23").Value = 0.3232848336488418
$ws.Range("J23").Value = 0.3232848336488418
$ws.Range("O23").Value = 0.1342711086924142
$ws.Range("P23").Value = 0.1342711086924142
$ws.Range("Q23").Value = 301.208340029226
$ws.Range("R23").Value = 2710.875060263034
$ws.Range("S23").Value = 0.04340781303747269
$ws.Range("T23").Value = 0.04340781303747267
$ws.Range("G24").Value = 41.78521799999999
$ws.Range("H24").Value = 125.355654
$ws.Range("I24").Value = 0.3232848336488418
$ws.Range("J24").Value = 0.3232848336488418
$ws.Range("M24").Value = 11.375406
$ws.Range("N24").Value = 34.126218
$ws.Range("O24").Value = 0.2118874139822907
$ws.Range("P24").Value = 0.2118874139822907
$ws.Range("Q24").Value = 475.3238195485079
$ws.Range("R24").Value = 4277.914375936572
$ws.Range("S24").Value = 0.06849998738154814
$ws.Range("T24").Value = 0.06849998738154811
$ws.Range("G25").Value = 41.78521799999999
$ws.Range("H25").Value = 125.355654
$ws.Range("I25").Value = 0.3232848336488418
$ws.Range("J25").Value = 0.3232848336488418
$ws.Range("M25").Value = 3.401340666666667
$ws.Range("N25").Value = 10.204022
$ws.Range("O25").Value = 0.06335609277882483
$ws.Range("P25").Value = 0.06335609277882483
$ws.Range("Q25").Value = 142.125761248932
$ws.Range("R25").Value = 1279.131851240388
$ws.Range("S25").Value = 0.02048206391464297
$ws.Range("T25").Value = 0.02048206391464297
$ws.Range("G26").Value = 41.78521799999999
$ws.Range("H26").Value = 125.355654
$ws.Range("I26").Value = 0.3232848336488418
$ws.Range("J26").Value = 0.3232848336488418
$ws.Range("M26").Value = 9.847764666666666
$ws.Range("N26").Value = 29.543294
$ws.Range("O26").Value = 0.1834323441929172
$ws.Range("P26").Value = 0.1834323441929172
$ws.Range("Q26").Value = 411.4909934093639
$ws.Range("R26").Value = 3703.418940684276
$ws.Range("S26").Value = 0.05930089487822431
$ws.Range("T26").Value = 0.05930089487822431
